# Generate Report for handback
# Updates the localization-status workbook to reflect that the two files
# have been handed back and are now in sync with en-US:
#  - Status text changes from "Ready for handoff" to
#    "Handed back: in sync with en-US" on every sheet that shows it.
#  - The zh-cn and de-de sheets gain "Latest Target File" (E) and
#    "Latest Handback File" (F) hyperlink entries for the two handed-off
#    files, mirroring the existing "Source File Name"/"Latest Handoff
#    File" hyperlinks.
#  - The "Latest Handback DateTime" column (G) is stamped with the
#    actual handback timestamp instead of the epoch placeholder.

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"

# Hyperlink colour used throughout the workbook for linked cells
# (matches the custom "HyperLink" cell style: underline, RGB 6495ED).
$linkColor = 15570276

function Set-LinkCell($ws, $cellAddr, $text, $url, $displayText) {
    $rng = $ws.Range($cellAddr)
    $rng.Value = $text
    $rng.Font.Underline = $true
    $rng.Font.Color = $linkColor
    $rng.Font.Name = "Calibri"
    $rng.Font.Size = 11
    $ws.Hyperlinks.Add($rng, $url, "", "", $displayText)
    $rng.Font.Underline = $true
    $rng.Font.Color = $linkColor
    $rng.Font.Name = "Calibri"
    $rng.Font.Size = 11
}

# ---------------------------------------------------------------------
# 1. Overview sheet - update Status columns (B, C) for rows 2 and 3
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = $statusText
$wsOverview.Range("C2").Value = $statusText
$wsOverview.Range("B3").Value = $statusText
$wsOverview.Range("C3").Value = $statusText

# ---------------------------------------------------------------------
# 2. zh-cn sheet
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("B2").Value = $statusText
$wsZh.Range("B3").Value = $statusText

Set-LinkCell $wsZh "E2" "72efcf84-c0da-4d37-b828-2cf8f4b5d32f.md" `
    "https://github.com/OpenLocalizationTest/oltest/blob/11ec0f254ad752e2561506e30edbb76110032dfd/e2e/72efcf84-c0da-4d37-b828-2cf8f4b5d32f.md" `
    "72efcf84-c0da-4d37-b828-2cf8f4b5d32f.md"

Set-LinkCell $wsZh "F2" "72efcf84-c0da-4d37-b828-2cf8f4b5d32f.a30c1bad72c67f13829bf83da99633040ce5a887.zh-cn.xlf" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c6ad097b06d7973a69eba93b3cecd8472df3cc60/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/72efcf84-c0da-4d37-b828-2cf8f4b5d32f.a30c1bad72c67f13829bf83da99633040ce5a887.zh-cn.xlf" `
    "72efcf84-c0da-4d37-b828-2cf8f4b5d32f.a30c1bad72c67f13829bf83da99633040ce5a887.zh-cn.xlf"

Set-LinkCell $wsZh "E3" "7783ce40-bf73-4bd0-b0ed-5c0dede8461b.md" `
    "https://github.com/OpenLocalizationTest/oltest/blob/11ec0f254ad752e2561506e30edbb76110032dfd/e2e/7783ce40-bf73-4bd0-b0ed-5c0dede8461b.md" `
    "7783ce40-bf73-4bd0-b0ed-5c0dede8461b.md"

Set-LinkCell $wsZh "F3" "7783ce40-bf73-4bd0-b0ed-5c0dede8461b.099770cfb3e9e332a48f0a66a7715163764e7df8.zh-cn.xlf" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c6ad097b06d7973a69eba93b3cecd8472df3cc60/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/7783ce40-bf73-4bd0-b0ed-5c0dede8461b.099770cfb3e9e332a48f0a66a7715163764e7df8.zh-cn.xlf" `
    "7783ce40-bf73-4bd0-b0ed-5c0dede8461b.099770cfb3e9e332a48f0a66a7715163764e7df8.zh-cn.xlf"

$wsZh.Range("G2").Value = "2016-02-16 13:53:23"
$wsZh.Range("G3").Value = "2016-02-16 13:53:23"

# ---------------------------------------------------------------------
# 3. de-de sheet
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("B2").Value = $statusText
$wsDe.Range("B3").Value = $statusText

Set-LinkCell $wsDe "E2" "72efcf84-c0da-4d37-b828-2cf8f4b5d32f.md" `
    "https://github.com/OpenLocalizationTest/oltest/blob/11ec0f254ad752e2561506e30edbb76110032dfd/e2e/72efcf84-c0da-4d37-b828-2cf8f4b5d32f.md" `
    "72efcf84-c0da-4d37-b828-2cf8f4b5d32f.md"

Set-LinkCell $wsDe "F2" "72efcf84-c0da-4d37-b828-2cf8f4b5d32f.a30c1bad72c67f13829bf83da99633040ce5a887.de-de.xlf" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/669ae8e51af2351bf57294f8461ab8dcb29f9064/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/72efcf84-c0da-4d37-b828-2cf8f4b5d32f.a30c1bad72c67f13829bf83da99633040ce5a887.de-de.xlf" `
    "72efcf84-c0da-4d37-b828-2cf8f4b5d32f.a30c1bad72c67f13829bf83da99633040ce5a887.de-de.xlf"

Set-LinkCell $wsDe "E3" "7783ce40-bf73-4bd0-b0ed-5c0dede8461b.md" `
    "https://github.com/OpenLocalizationTest/oltest/blob/11ec0f254ad752e2561506e30edbb76110032dfd/e2e/7783ce40-bf73-4bd0-b0ed-5c0dede8461b.md" `
    "7783ce40-bf73-4bd0-b0ed-5c0dede8461b.md"

Set-LinkCell $wsDe "F3" "7783ce40-bf73-4bd0-b0ed-5c0dede8461b.099770cfb3e9e332a48f0a66a7715163764e7df8.de-de.xlf" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/669ae8e51af2351bf57294f8461ab8dcb29f9064/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/7783ce40-bf73-4bd0-b0ed-5c0dede8461b.099770cfb3e9e332a48f0a66a7715163764e7df8.de-de.xlf" `
    "7783ce40-bf73-4bd0-b0ed-5c0dede8461b.099770cfb3e9e332a48f0a66a7715163764e7df8.de-de.xlf"

$wsDe.Range("G2").Value = "2016-02-16 13:53:51"
$wsDe.Range("G3").Value = "2016-02-16 13:53:51"
